$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AT2").Value = "['Netherlands', 1, 2]"
$ws.Range("AU2").Value = "['Austria', 3, 2]"

# Row 6
$ws.Range("AN6").Value = "['Italy', 1, 2]"
$ws.Range("AO6").Value = "['Germany', 2, 0]"
$ws.Range("AP6").Value = "['England', 2, 0]"
$ws.Range("AQ6").Value = "['Spain', 3, 0]"
$ws.Range("AR6").Value = "['France', 2, 1]"
$ws.Range("AS6").Value = "['Draw', 1, 1]"
$ws.Range("AT6").Value = "['Netherlands', 0, 1]"
$ws.Range("AU6").Value = "['Austria', 2, 1]"

# Row 7
$ws.Range("AN7").Value = "['Italy', 0, 2]"
$ws.Range("AO7").Value = "['Germany', 3, 0]"
$ws.Range("AP7").Value = "['England', 2, 0]"
$ws.Range("AQ7").Value = "['Spain', 3, 1]"
$ws.Range("AR7").Value = "['France', 2, 1]"
$ws.Range("AS7").Value = "['Portugal', 2, 0]"
$ws.Range("AT7").Value = "['Netherlands', 0, 1]"
$ws.Range("AU7").Value = "['Austria', 2, 1]"

# Row 8
$ws.Range("AN8").Value = "['Italy', 0, 2]"
$ws.Range("AO8").Value = "['Germany', 1, 0]"
$ws.Range("AP8").Value = "['England', 1, 0]"
$ws.Range("AQ8").Value = "['Spain', 3, 1]"
$ws.Range("AR8").Value = "['Draw', 2, 2]"
$ws.Range("AS8").Value = "['Portugal', 3, 1]"
$ws.Range("AT8").Value = "['Netherlands', 0, 2]"
$ws.Range("AU8").Value = "['Draw', 1, 1]"

# Row 10
$ws.Range("AN10").Value = "['Italy', 1, 2]"
$ws.Range("AO10").Value = "['Germany', 3, 0]"
$ws.Range("AP10").Value = "['Draw', 0, 0]"
$ws.Range("AQ10").Value = "['Spain', 2, 1]"
$ws.Range("AR10").Value = "['Draw', 1, 1]"
$ws.Range("AS10").Value = "['Portugal', 2, 0]"
$ws.Range("AT10").Value = "['Netherlands', 0, 2]"
$ws.Range("AU10").Value = "['Austria', 2, 1]"

# Row 11
$ws.Range("AN11").Value = "['Draw', 1, 1]"
$ws.Range("AO11").Value = "['Germany', 2, 1]"
$ws.Range("AP11").Value = "['England', 1, 0]"
$ws.Range("AQ11").Value = "['Spain', 2, 0]"
$ws.Range("AR11").Value = "['France', 2, 0]"
$ws.Range("AS11").Value = "['Portugal', 1, 0]"
$ws.Range("AT11").Value = "['Netherlands', 1, 2]"
$ws.Range("AU11").Value = "['Austria', 2, 1]"

# Row 13 (partial: AR13 and AT13 stay empty)
$ws.Range("AN13").Value = "['Draw', 0, 0]"
$ws.Range("AO13").Value = "['Germany', 2, 0]"
$ws.Range("AP13").Value = "['England', 1, 0]"
$ws.Range("AQ13").Value = "['Spain', 3, 1]"
$ws.Range("AS13").Value = "['Portugal', 2, 0]"
$ws.Range("AU13").Value = "['Austria', 3, 2]"

# Row 19 (partial: only AR19-AU19 set)
$ws.Range("AR19").Value = "['France', 3, 0]"
$ws.Range("AS19").Value = "['Portugal', 2, 1]"
$ws.Range("AT19").Value = "['Draw', 1, 1]"
$ws.Range("AU19").Value = "['Austria', 1, 0]"

# Row 21
$ws.Range("AN21").Value = "['Italy', 1, 2]"
$ws.Range("AO21").Value = "['Germany', 2, 1]"
$ws.Range("AP21").Value = "['England', 2, 0]"
$ws.Range("AQ21").Value = "['Spain', 3, 1]"
$ws.Range("AR21").Value = "['France', 3, 2]"
$ws.Range("AS21").Value = "['Portugal', 2, 1]"
$ws.Range("AT21").Value = "['Netherlands', 1, 2]"
$ws.Range("AU21").Value = "['Austria', 3, 1]"

# Row 23
$ws.Range("AN23").Value = "['Italy', 0, 1]"
$ws.Range("AO23").Value = "['Germany', 2, 1]"
$ws.Range("AP23").Value = "['England', 2, 0]"
$ws.Range("AQ23").Value = "['Spain', 2, 0]"
$ws.Range("AR23").Value = "['Draw', 1, 1]"
$ws.Range("AS23").Value = "['Portugal', 2, 0]"
$ws.Range("AT23").Value = "['Netherlands', 0, 1]"
$ws.Range("AU23").Value = "['Austria', 2, 1]"

# Row 30 (modify existing values)
$ws.Range("AS30").Value = "['Portugal', 2, 0]"
$ws.Range("AT30").Value = "['Netherlands', 1, 2]"

# Row 31 (partial: AN31/AO31 unchanged, AP31-AU31 newly set)
$ws.Range("AP31").Value = "['England', 3, 0]"
$ws.Range("AQ31").Value = "['Spain', 3, 1]"
$ws.Range("AR31").Value = "['France', 1, 0]"
$ws.Range("AS31").Value = "['Portugal', 2, 0]"
$ws.Range("AT31").Value = "['Netherlands', 0, 1]"
$ws.Range("AU31").Value = "['Austria', 2, 1]"
